$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 30.75612566666667
$ws.Cells.Item(2, 8).Value = 92.268377
$ws.Cells.Item(2, 9).Value = 0.9777985798685588
$ws.Cells.Item(2, 10).Value = 0.9777985798685588
$ws.Cells.Item(2, 13).Value = 68.63737500000001
$ws.Cells.Item(2, 14).Value = 205.912125
$ws.Cells.Item(2, 15).Value = 0.5415701538216162
$ws.Cells.Item(2, 16).Value = 0.5415701538216162
$ws.Cells.Item(2, 17).Value = 2111.019730930125
$ws.Cells.Item(2, 18).Value = 18999.17757837113
$ws.Cells.Item(2, 19).Value = 0.5295465273059733
$ws.Cells.Item(2, 20).Value = 0.5295465273059733

# Row 3
$ws.Cells.Item(3, 7).Value = 30.75612566666667
$ws.Cells.Item(3, 8).Value = 92.268377
$ws.Cells.Item(3, 9).Value = 0.9777985798685588
$ws.Cells.Item(3, 10).Value = 0.9777985798685588
$ws.Cells.Item(3, 15).Value = 0.08718851262838957
$ws.Cells.Item(3, 16).Value = 0.08718851262838957
$ws.Cells.Item(3, 17).Value = 339.8574850740498
$ws.Cells.Item(3, 18).Value = 3058.717365666448
$ws.Cells.Item(3, 19).Value = 0.08525280382889122
$ws.Cells.Item(3, 20).Value = 0.08525280382889122

# Row 4
$ws.Cells.Item(4, 7).Value = 30.75612566666667
$ws.Cells.Item(4, 8).Value = 92.268377
$ws.Cells.Item(4, 9).Value = 0.9777985798685588
$ws.Cells.Item(4, 10).Value = 0.9777985798685588
$ws.Cells.Item(4, 13).Value = 16.21089566666667
$ws.Cells.Item(4, 14).Value = 48.632687
$ws.Cells.Item(4, 15).Value = 0.1279089892319285
$ws.Cells.Item(4, 16).Value = 0.1279089892319285
$ws.Cells.Item(4, 17).Value = 498.5843442932222
$ws.Cells.Item(4, 18).Value = 4487.259098639
$ws.Cells.Item(4, 19).Value = 0.1250692280234025
$ws.Cells.Item(4, 20).Value = 0.1250692280234025

# Row 5
$ws.Cells.Item(5, 7).Value = 30.75612566666667
$ws.Cells.Item(5, 8).Value = 92.268377
$ws.Cells.Item(5, 9).Value = 0.9777985798685588
$ws.Cells.Item(5, 10).Value = 0.9777985798685588
$ws.Cells.Item(5, 13).Value = 20.32546233333333
$ws.Cells.Item(5, 14).Value = 60.976387
$ws.Cells.Item(5, 15).Value = 0.1603741949973873
$ws.Cells.Item(5, 16).Value = 0.1603741949973873
$ws.Cells.Item(5, 17).Value = 625.1324737570999
$ws.Cells.Item(5, 18).Value = 5626.192263813899
$ws.Cells.Item(5, 19).Value = 0.1568136601160087
$ws.Cells.Item(5, 20).Value = 0.1568136601160087

# Row 6
$ws.Cells.Item(6, 7).Value = 30.75612566666667
$ws.Cells.Item(6, 8).Value = 92.268377
$ws.Cells.Item(6, 9).Value = 0.9777985798685588
$ws.Cells.Item(6, 10).Value = 0.9777985798685588
$ws.Cells.Item(6, 13).Value = 10.513928
$ws.Cells.Item(6, 14).Value = 31.541784
$ws.Cells.Item(6, 15).Value = 0.08295814932067838
$ws.Cells.Item(6, 16).Value = 0.08295814932067838
$ws.Cells.Item(6, 17).Value = 323.3676908182853
$ws.Cells.Item(6, 18).Value = 2910.309217364568
$ws.Cells.Item(6, 19).Value = 0.08111636059428316
$ws.Cells.Item(6, 20).Value = 0.08111636059428316

# Row 7
$ws.Cells.Item(7, 9).Value = 0.004830327290741966
$ws.Cells.Item(7, 10).Value = 0.004830327290741966
$ws.Cells.Item(7, 13).Value = 68.63737500000001
$ws.Cells.Item(7, 14).Value = 205.912125
$ws.Cells.Item(7, 15).Value = 0.5415701538216162
$ws.Cells.Item(7, 16).Value = 0.5415701538216162
$ws.Cells.Item(7, 17).Value = 10.42844244975
$ws.Cells.Item(7, 18).Value = 93.85598204775
$ws.Cells.Item(7, 19).Value = 0.002615961093855877
$ws.Cells.Item(7, 20).Value = 0.002615961093855877

# Row 8
$ws.Cells.Item(8, 9).Value = 0.004830327290741966
$ws.Cells.Item(8, 10).Value = 0.004830327290741966
$ws.Cells.Item(8, 15).Value = 0.08718851262838957
$ws.Cells.Item(8, 16).Value = 0.08718851262838957
$ws.Cells.Item(8, 19).Value = 0.0004211490519881107
$ws.Cells.Item(8, 20).Value = 0.0004211490519881107

# Row 9
$ws.Cells.Item(9, 9).Value = 0.004830327290741966
$ws.Cells.Item(9, 10).Value = 0.004830327290741966
$ws.Cells.Item(9, 13).Value = 16.21089566666667
$ws.Cells.Item(9, 14).Value = 48.632687
$ws.Cells.Item(9, 15).Value = 0.1279089892319285
$ws.Cells.Item(9, 16).Value = 0.1279089892319285
$ws.Cells.Item(9, 17).Value = 2.463007836746889
$ws.Cells.Item(9, 18).Value = 22.167070530722
$ws.Cells.Item(9, 19).Value = 0.0006178422814182045
$ws.Cells.Item(9, 20).Value = 0.0006178422814182045

# Row 10
$ws.Cells.Item(10, 9).Value = 0.004830327290741966
$ws.Cells.Item(10, 10).Value = 0.004830327290741966
$ws.Cells.Item(10, 13).Value = 20.32546233333333
$ws.Cells.Item(10, 14).Value = 60.976387
$ws.Cells.Item(10, 15).Value = 0.1603741949973873
$ws.Cells.Item(10, 16).Value = 0.1603741949973873
$ws.Cells.Item(10, 17).Value = 3.088155894769111
$ws.Cells.Item(10, 18).Value = 27.793403052922
$ws.Cells.Item(10, 19).Value = 0.0007746598508266537
$ws.Cells.Item(10, 20).Value = 0.0007746598508266537

# Row 11
$ws.Cells.Item(11, 9).Value = 0.004830327290741966
$ws.Cells.Item(11, 10).Value = 0.004830327290741966
$ws.Cells.Item(11, 13).Value = 10.513928
$ws.Cells.Item(11, 14).Value = 31.541784
$ws.Cells.Item(11, 15).Value = 0.08295814932067838
$ws.Cells.Item(11, 16).Value = 0.08295814932067838
$ws.Cells.Item(11, 17).Value = 1.597437155322667
$ws.Cells.Item(11, 18).Value = 14.376934397904
$ws.Cells.Item(11, 19).Value = 0.0004007150126531199
$ws.Cells.Item(11, 20).Value = 0.0004007150126531199

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.4673666666666667
$ws.Cells.Item(12, 8).Value = 1.4021
$ws.Cells.Item(12, 9).Value = 0.01485851852399773
$ws.Cells.Item(12, 10).Value = 0.01485851852399773
$ws.Cells.Item(12, 13).Value = 68.63737500000001
$ws.Cells.Item(12, 14).Value = 205.912125
$ws.Cells.Item(12, 15).Value = 0.5415701538216162
$ws.Cells.Item(12, 16).Value = 0.5415701538216162
$ws.Cells.Item(12, 17).Value = 32.07882116250001
$ws.Cells.Item(12, 18).Value = 288.7093904625
$ws.Cells.Item(12, 19).Value = 0.008046930162602786
$ws.Cells.Item(12, 20).Value = 0.008046930162602786

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.4673666666666667
$ws.Cells.Item(13, 8).Value = 1.4021
$ws.Cells.Item(13, 9).Value = 0.01485851852399773
$ws.Cells.Item(13, 10).Value = 0.01485851852399773
$ws.Cells.Item(13, 15).Value = 0.08718851262838957
$ws.Cells.Item(13, 16).Value = 0.08718851262838957
$ws.Cells.Item(13, 17).Value = 5.164436563377778
$ws.Cells.Item(13, 18).Value = 46.4799290704
$ws.Cells.Item(13, 19).Value = 0.001295492129968737
$ws.Cells.Item(13, 20).Value = 0.001295492129968737

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.4673666666666667
$ws.Cells.Item(14, 8).Value = 1.4021
$ws.Cells.Item(14, 9).Value = 0.01485851852399773
$ws.Cells.Item(14, 10).Value = 0.01485851852399773
$ws.Cells.Item(14, 13).Value = 16.21089566666667
$ws.Cells.Item(14, 14).Value = 48.632687
$ws.Cells.Item(14, 15).Value = 0.1279089892319285
$ws.Cells.Item(14, 16).Value = 0.1279089892319285
$ws.Cells.Item(14, 17).Value = 7.576432271411113
$ws.Cells.Item(14, 18).Value = 68.18789044270001
$ws.Cells.Item(14, 19).Value = 0.001900538085888436
$ws.Cells.Item(14, 20).Value = 0.001900538085888436

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.4673666666666667
$ws.Cells.Item(15, 8).Value = 1.4021
$ws.Cells.Item(15, 9).Value = 0.01485851852399773
$ws.Cells.Item(15, 10).Value = 0.01485851852399773
$ws.Cells.Item(15, 13).Value = 20.32546233333333
$ws.Cells.Item(15, 14).Value = 60.976387
$ws.Cells.Item(15, 15).Value = 0.1603741949973873
$ws.Cells.Item(15, 16).Value = 0.1603741949973873
$ws.Cells.Item(15, 17).Value = 9.49944357918889
$ws.Cells.Item(15, 18).Value = 85.49499221270001
$ws.Cells.Item(15, 19).Value = 0.002382922947139904
$ws.Cells.Item(15, 20).Value = 0.002382922947139904

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.4673666666666667
$ws.Cells.Item(16, 8).Value = 1.4021
$ws.Cells.Item(16, 9).Value = 0.01485851852399773
$ws.Cells.Item(16, 10).Value = 0.01485851852399773
$ws.Cells.Item(16, 13).Value = 10.513928
$ws.Cells.Item(16, 14).Value = 31.541784
$ws.Cells.Item(16, 15).Value = 0.08295814932067838
$ws.Cells.Item(16, 16).Value = 0.08295814932067838
$ws.Cells.Item(16, 17).Value = 4.913859482933334
$ws.Cells.Item(16, 18).Value = 44.2247353464
$ws.Cells.Item(16, 19).Value = 0.00123263519839787
$ws.Cells.Item(16, 20).Value = 0.00123263519839787

# Row 17
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.07903166666666667
$ws.Cells.Item(17, 8).Value = 0.237095
$ws.Cells.Item(17, 9).Value = 0.002512574316701549
$ws.Cells.Item(17, 10).Value = 0.002512574316701549
$ws.Cells.Item(17, 13).Value = 68.63737500000001
$ws.Cells.Item(17, 14).Value = 205.912125
$ws.Cells.Item(17, 15).Value = 0.5415701538216162
$ws.Cells.Item(17, 16).Value = 0.5415701538216162
$ws.Cells.Item(17, 17).Value = 5.424526141875001
$ws.Cells.Item(17, 18).Value = 48.820735276875
$ws.Cells.Item(17, 19).Value = 0.0013607352591843
$ws.Cells.Item(17, 20).Value = 0.0013607352591843

# Row 18
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 0.6666666666666666
$ws.Cells.Item(18, 7).Value = 0.07903166666666667
$ws.Cells.Item(18, 8).Value = 0.237095
$ws.Cells.Item(18, 9).Value = 0.002512574316701549
$ws.Cells.Item(18, 10).Value = 0.002512574316701549
$ws.Cells.Item(18, 15).Value = 0.08718851262838957
$ws.Cells.Item(18, 16).Value = 0.08718851262838957
$ws.Cells.Item(18, 17).Value = 0.8733058176977778
$ws.Cells.Item(18, 18).Value = 7.859752359280001
$ws.Cells.Item(18, 19).Value = 0.0002190676175415003
$ws.Cells.Item(18, 20).Value = 0.0002190676175415003

# Row 19
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 0.6666666666666666
$ws.Cells.Item(19, 7).Value = 0.07903166666666667
$ws.Cells.Item(19, 8).Value = 0.237095
$ws.Cells.Item(19, 9).Value = 0.002512574316701549
$ws.Cells.Item(19, 10).Value = 0.002512574316701549
$ws.Cells.Item(19, 13).Value = 16.21089566666667
$ws.Cells.Item(19, 14).Value = 48.632687
$ws.Cells.Item(19, 15).Value = 0.1279089892319285
$ws.Cells.Item(19, 16).Value = 0.1279089892319285
$ws.Cells.Item(19, 17).Value = 1.281174102696111
$ws.Cells.Item(19, 18).Value = 11.530566924265
$ws.Cells.Item(19, 19).Value = 0.0003213808412193986
$ws.Cells.Item(19, 20).Value = 0.0003213808412193986

# Row 20
$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = 0.6666666666666666
$ws.Cells.Item(20, 7).Value = 0.07903166666666667
$ws.Cells.Item(20, 8).Value = 0.237095
$ws.Cells.Item(20, 9).Value = 0.002512574316701549
$ws.Cells.Item(20, 10).Value = 0.002512574316701549
$ws.Cells.Item(20, 13).Value = 20.32546233333333
$ws.Cells.Item(20, 14).Value = 60.976387
$ws.Cells.Item(20, 15).Value = 0.1603741949973873
$ws.Cells.Item(20, 16).Value = 0.1603741949973873
$ws.Cells.Item(20, 17).Value = 1.606355163973889
$ws.Cells.Item(20, 18).Value = 14.457196475765
$ws.Cells.Item(20, 19).Value = 0.0004029520834121215
$ws.Cells.Item(20, 20).Value = 0.0004029520834121215

# Row 21
$ws.Cells.Item(21, 5).Value = 2
$ws.Cells.Item(21, 6).Value = 0.6666666666666666
$ws.Cells.Item(21, 7).Value = 0.07903166666666667
$ws.Cells.Item(21, 8).Value = 0.237095
$ws.Cells.Item(21, 9).Value = 0.002512574316701549
$ws.Cells.Item(21, 10).Value = 0.002512574316701549
$ws.Cells.Item(21, 13).Value = 10.513928
$ws.Cells.Item(21, 14).Value = 31.541784
$ws.Cells.Item(21, 15).Value = 0.08295814932067838
$ws.Cells.Item(21, 16).Value = 0.08295814932067838
$ws.Cells.Item(21, 17).Value = 0.8309332530533333
$ws.Cells.Item(21, 18).Value = 7.47839927748
$ws.Cells.Item(21, 19).Value = 0.0002084385153442286
$ws.Cells.Item(21, 20).Value = 0.0002084385153442286
